$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log entries (rows shift down; new rows 3-6 inserted above the old row3 "Kullu..." entry) ---
# Cells are written in the same order the original author's Excel session appended them to the
# shared-string table (bottom row to top row), so the resulting sharedStrings.xml ordering matches.

# Row 6: 2021-12-12 normalization/imputation
$ws.Range("F6").Value = "normalization = {float, reflectance, normalize, standardize}; imputation = {fill 0, forward filling, linear interpolation}"
$ws.Range("C6").Value = "cropland_1206-215716.log -- cropland_1212-232444.log"
$ws.Range("B6").Value = "fa718209df583a9c48b2c825f25b310ef0468816"
$ws.Range("E6").Value = "Try different normalization and imputation methods"
$ws.Range("A6").Value = 44542
$ws.Range("D6").Value = "tile 43SFR"

# Row 5: 2022-01-04 cloud mask
$ws.Range("B5").Value = "54d2e4be4fbad9f623166c0f80e8b565da3df7ab"
$ws.Range("C5").Value = "cropland_20220104-112305.log"
$ws.Range("E5").Value = "Change cloud mask"
$ws.Range("F5").Value = "use scene classification results from Sen2Cor rather than the built-in mask"
$ws.Range("A5").Value = 44565
$ws.Range("D5").Value = "tile 43SFR"

# Row 4: 2022-01-11 feature selection
$ws.Range("E4").Value = "Feature selection"
$ws.Range("F4").Value = "run feature selection for rfc: {temporal, temporal+ndvi_spatial, temporal+spatial, select}"
$ws.Range("C4").Value = "cropland_20220105-135132.log -- cropland_20220111-094108.log"
$ws.Range("B4").Value = "139c94e3e54f8999023c00cdd34e4e92552df70b"
$ws.Range("A4").Value = 44572
$ws.Range("D4").Value = "tile 43SFR"

# Row 3: 2022-01-17 predict on 4 tiles
$ws.Range("C3").Value = "cropland_20220117-172418_predict.log"
$ws.Range("E3").Value = "Generate cropland map"
$ws.Range("F3").Value = "make cropland predictions on 4 tiles"
$ws.Range("D3").Value = "tiles 43SFR, 43RFQ, 43SGR, 43SGR"
$ws.Range("B3").Value = "af9bf4de1aa966d1ce515c6a727f2238581c5c7a"
$ws.Range("A3").Value = 44578

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 17.5
$ws.Rows.Item(3).RowHeight = 43.5
$ws.Rows.Item(4).RowHeight = 43.5
$ws.Rows.Item(5).RowHeight = 43.5
$ws.Rows.Item(6).RowHeight = 43.5
$ws.Rows.Item(7).RowHeight = 58
$ws.Rows.Item(8).RowHeight = 58
$ws.Rows.Item(9).RowHeight = 43.5
$ws.Rows.Item(10).RowHeight = 43.5
$ws.Rows.Item(11).RowHeight = 43.5

# --- Selection / view state ---
$ws.Range("C3").Select()
